# Updated cryptos list on Fri Jun  9 17:41:16 UTC 2023 with GitHub Actions
# Refreshes the per-coin Price (D) / Volume(1h) (E) figures, and for the two
# rows whose rank order flipped this run (EnergySwap/Decentraland and
# Elrond/EOS) also rewrites the Coin (B) and Link (C) cells so each row's
# data stays together while column A's rank index is left untouched.
#
# Cells whose new Price text would otherwise be auto-parsed by Excel as a
# number (e.g. "9.089") are explicitly formatted as Text first so the
# stored value keeps matching the original plain-text "26.523.37"-style
# price strings used throughout column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.523.37'
$ws.Range('E2').Value = '  +0.08%  '

$ws.Range('D3').Value = '1.843.54'
$ws.Range('E3').Value = '  -0.14%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '262.43'
$ws.Range('E5').Value = '  -0.34%  '

$ws.Range('E6').Value = '  +0.01%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5327'
$ws.Range('E7').Value = '  +2.31%  '

$ws.Range('E8').Value = '  -4.88%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06894'
$ws.Range('E9').Value = '  +1.78%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.41'
$ws.Range('E10').Value = '  -1.22%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07811'
$ws.Range('E11').Value = '  +0.43%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7540'
$ws.Range('E12').Value = '  -1.99%  '

$ws.Range('D13').Value = '1.844.27'
$ws.Range('E13').Value = '  -0.85%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '89.90'
$ws.Range('E14').Value = '  +1.78%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.016'
$ws.Range('E15').Value = '  +0.10%  '

$ws.Range('E16').Value = '  -0.05%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.01'
$ws.Range('E17').Value = '  +0.63%  '

$ws.Range('E18').Value = '  +0.14%  '

$ws.Range('D20').Value = '26.547.76'
$ws.Range('E20').Value = '  +0.01%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.619'
$ws.Range('E21').Value = '  +0.17%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.989'
$ws.Range('E22').Value = '  +0.28%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.330'
$ws.Range('E23').Value = '  -1.05%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '142.95'
$ws.Range('E24').Value = '  -0.23%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.193'
$ws.Range('E25').Value = '  +0.48%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.684'
$ws.Range('E26').Value = '  +0.37%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.97'
$ws.Range('E27').Value = '  +0.11%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '110.85'
$ws.Range('E28').Value = '  -0.88%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.282'
$ws.Range('E29').Value = '  +2.96%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.08808'
$ws.Range('E30').Value = '  +0.88%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.088'
$ws.Range('E31').Value = '  -0.40%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.04813'
$ws.Range('E32').Value = '  +0.04%  '

$ws.Range('E33').Value = '  +2.23%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7314'
$ws.Range('E34').Value = '  +2.38%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.135'
$ws.Range('E35').Value = '  +0.64%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.108'
$ws.Range('E36').Value = '  +0.38%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.304'
$ws.Range('E37').Value = '  +5.13%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01718'
$ws.Range('E38').Value = '  -3.60%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.4796'

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9086'
$ws.Range('E40').Value = '  +1.44%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '108.42'
$ws.Range('E41').Value = '  -3.44%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.870'
$ws.Range('E42').Value = '  -2.55%  '

$ws.Range('E43').Value = '  +0.01%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '7.509'
$ws.Range('E44').Value = '  -1.23%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '9.089'
$ws.Range('E45').Value = '  +0.45%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4143'
$ws.Range('E46').Value = '  -0.44%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.1244'
$ws.Range('E47').Value = '  +1.48%  '

$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '34.89'
$ws.Range('E48').Value = '  -0.06%  '

$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.9001'
$ws.Range('E49').Value = '  +1.79%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05790'
$ws.Range('E50').Value = '  -1.87%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '60.33'
$ws.Range('E51').Value = '  +0.91%  '
